{"js": "const body = context.document.body;\n\n{\n  const results = body.search(\"781\u00d73=2343\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"391\u00d73=1173\", Word.InsertLocation.replace);\n  }\n}\n{\n  const results = body.search(\"911\u00d72=1822\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"816\u00d77=5712\", Word.InsertLocation.replace);\n  }\n}\n{\n  const results = body.search(\"959\u00d76=5754\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"874\u00d77=6118\", Word.InsertLocation.replace);\n  }\n}\n{\n  const results = body.search(\"963\u00d74=3852\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"936\u00d75=4680\", Word.InsertLocation.replace);\n  }\n}\n{\n  const results = body.search(\"206\u00d72=412\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"177\u00d78=1416\", Word.InsertLocation.replace);\n  }\n}\n{\n  const results = body.search(\"166\u00d76=996\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"173\u00d72=346\", Word.InsertLocation.replace);\n  }\n}\n{\n  const results = body.search(\"230\u00d74=920\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"744\u00d77=5208\", Word.InsertLocation.replace);\n  }\n}\n{\n  const results = body.search(\"756\u00d76=4536\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"770\u00d77=5390\", Word.InsertLocation.replace);\n  }\n}\n{\n  const results = body.search(\"966\u00d72=1932\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"971\u00d78=7768\", Word.InsertLocation.replace);\n  }\n}\n{\n  const results = body.search(\"341\u00d72=682\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"390\u00d73=1170\", Word.InsertLocation.replace);\n  }\n}\n{\n  const results = body.search(\"201\u00d77=1407\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"551\u00d76=3306\", Word.InsertLocation.replace);\n  }\n}\n{\n  const results = body.search(\"803\u00d77=5621\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"917\u00d74=3668\", Word.InsertLocation.replace);\n  }\n}\n{\n  const results = body.search(\"543\u00d79=4887\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"248\u00d78=1984\", Word.InsertLocation.replace);\n  }\n}\n{\n  const results = body.search(\"275\u00d72=550\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"248\u00d77=1736\", Word.InsertLocation.replace);\n  }\n}\n{\n  const results = body.search(\"287\u00d78=2296\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"826\u00d73=2478\", Word.InsertLocation.replace);\n  }\n}\n{\n  const results = body.search(\"693\u00d78=5544\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"486\u00d77=3402\", Word.InsertLocation.replace);\n  }\n}\n{\n  const results = body.search(\"585\u00d74=2340\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"910\u00d74=3640\", Word.InsertLocation.replace);\n  }\n}\n{\n  const results = body.search(\"932\u00d76=5592\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"645\u00d73=1935\", Word.InsertLocation.replace);\n  }\n}\n{\n  const results = body.search(\"154\u00d78=1232\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"627\u00d78=5016\", Word.InsertLocation.replace);\n  }\n}\n{\n  const results = body.search(\"522\u00d72=1044\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"394\u00d73=1182\", Word.InsertLocation.replace);\n  }\n}\n{\n  const results = body.search(\"900\u00d75=4500\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"307\u00d74=1228\", Word.InsertLocation.replace);\n  }\n}\n{\n  const results = body.search(\"478\u00d76=2868\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"668\u00d78=5344\", Word.InsertLocation.replace);\n  }\n}\n{\n  const results = body.search(\"301\u00d78=2408\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"365\u00d74=1460\", Word.InsertLocation.replace);\n  }\n}\n{\n  const results = body.search(\"382\u00d76=2292\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"497\u00d79=4473\", Word.InsertLocation.replace);\n  }\n}\n{\n  const results = body.search(\"364\u00d74=1456\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"925\u00d78=7400\", Word.InsertLocation.replace);\n  }\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$d.Content.Find.Execute(\"781\u00d73=2343\", $false, $false, $false, $false, $false, $true, 1, $false, \"391\u00d73=1173\", 2) | Out-Null\n$d.Content.Find.Execute(\"911\u00d72=1822\", $false, $false, $false, $false, $false, $true, 1, $false, \"816\u00d77=5712\", 2) | Out-Null\n$d.Content.Find.Execute(\"959\u00d76=5754\", $false, $false, $false, $false, $false, $true, 1, $false, \"874\u00d77=6118\", 2) | Out-Null\n$d.Content.Find.Execute(\"963\u00d74=3852\", $false, $false, $false, $false, $false, $true, 1, $false, \"936\u00d75=4680\", 2) | Out-Null\n$d.Content.Find.Execute(\"206\u00d72=412\", $false, $false, $false, $false, $false, $true, 1, $false, \"177\u00d78=1416\", 2) | Out-Null\n$d.Content.Find.Execute(\"166\u00d76=996\", $false, $false, $false, $false, $false, $true, 1, $false, \"173\u00d72=346\", 2) | Out-Null\n$d.Content.Find.Execute(\"230\u00d74=920\", $false, $false, $false, $false, $false, $true, 1, $false, \"744\u00d77=5208\", 2) | Out-Null\n$d.Content.Find.Execute(\"756\u00d76=4536\", $false, $false, $false, $false, $false, $true, 1, $false, \"770\u00d77=5390\", 2) | Out-Null\n$d.Content.Find.Execute(\"966\u00d72=1932\", $false, $false, $false, $false, $false, $true, 1, $false, \"971\u00d78=7768\", 2) | Out-Null\n$d.Content.Find.Execute(\"341\u00d72=682\", $false, $false, $false, $false, $false, $true, 1, $false, \"390\u00d73=1170\", 2) | Out-Null\n$d.Content.Find.Execute(\"201\u00d77=1407\", $false, $false, $false, $false, $false, $true, 1, $false, \"551\u00d76=3306\", 2) | Out-Null\n$d.Content.Find.Execute(\"803\u00d77=5621\", $false, $false, $false, $false, $false, $true, 1, $false, \"917\u00d74=3668\", 2) | Out-Null\n$d.Content.Find.Execute(\"543\u00d79=4887\", $false, $false, $false, $false, $false, $true, 1, $false, \"248\u00d78=1984\", 2) | Out-Null\n$d.Content.Find.Execute(\"275\u00d72=550\", $false, $false, $false, $false, $false, $true, 1, $false, \"248\u00d77=1736\", 2) | Out-Null\n$d.Content.Find.Execute(\"287\u00d78=2296\", $false, $false, $false, $false, $false, $true, 1, $false, \"826\u00d73=2478\", 2) | Out-Null\n$d.Content.Find.Execute(\"693\u00d78=5544\", $false, $false, $false, $false, $false, $true, 1, $false, \"486\u00d77=3402\", 2) | Out-Null\n$d.Content.Find.Execute(\"585\u00d74=2340\", $false, $false, $false, $false, $false, $true, 1, $false, \"910\u00d74=3640\", 2) | Out-Null\n$d.Content.Find.Execute(\"932\u00d76=5592\", $false, $false, $false, $false, $false, $true, 1, $false, \"645\u00d73=1935\", 2) | Out-Null\n$d.Content.Find.Execute(\"154\u00d78=1232\", $false, $false, $false, $false, $false, $true, 1, $false, \"627\u00d78=5016\", 2) | Out-Null\n$d.Content.Find.Execute(\"522\u00d72=1044\", $false, $false, $false, $false, $false, $true, 1, $false, \"394\u00d73=1182\", 2) | Out-Null\n$d.Content.Find.Execute(\"900\u00d75=4500\", $false, $false, $false, $false, $false, $true, 1, $false, \"307\u00d74=1228\", 2) | Out-Null\n$d.Content.Find.Execute(\"478\u00d76=2868\", $false, $false, $false, $false, $false, $true, 1, $false, \"668\u00d78=5344\", 2) | Out-Null\n$d.Content.Find.Execute(\"301\u00d78=2408\", $false, $false, $false, $false, $false, $true, 1, $false, \"365\u00d74=1460\", 2) | Out-Null\n$d.Content.Find.Execute(\"382\u00d76=2292\", $false, $false, $false, $false, $false, $true, 1, $false, \"497\u00d79=4473\", 2) | Out-Null\n$d.Content.Find.Execute(\"364\u00d74=1456\", $false, $false, $false, $false, $false, $true, 1, $false, \"925\u00d78=7400\", 2) | Out-Null\n\nWrite-Output \"Replaced 25 multiplication facts\"\n"}
